$wb = $excel.ActiveWorkbook

# The real edit duplicated the "week 4" sheet to start a new week's log,
# renamed the copy "week 7", placed it right before "Totaal", filled in the
# first row with the new entry and cleared the rest of the inherited sample
# data, then left "week 7" as the active/selected sheet.
#
# Reserve sheetId 8 with a throwaway sheet first so that when we copy
# "week 4" right after, the new copy receives sheetId 9 (matching the
# ever-increasing internal sheetId counter Excel itself uses), then remove
# the throwaway sheet.
$dummy = $wb.Worksheets.Add()

$src = $wb.Worksheets.Item("week 4")
$totaal = $wb.Worksheets.Item("Totaal")
$src.Copy($totaal)

$dummy.Delete()

$ws = $wb.Worksheets.Item("week 4 (2)")
$ws.Name = "week 7"

# Fill in the new week's first log entry (row 7).
$ws.Range("B7").Value2 = 41681
$ws.Range("C7").Value2 = 0.36458333333333331
$ws.Range("F7").Value2 = "bezig in upload_form"

# Clear the rest of the copied-over sample data from "week 4" (rows 8, 10-14)
# while leaving the row/column structure, styles and helper formulas intact.
$ws.Range("C8").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("F8").ClearContents()

$ws.Range("B10").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("D10").ClearContents()
$ws.Range("F10").ClearContents()

$ws.Range("C11").ClearContents()
$ws.Range("D11").ClearContents()
$ws.Range("F11").ClearContents()

$ws.Range("C12").ClearContents()
$ws.Range("D12").ClearContents()
$ws.Range("F12").ClearContents()

$ws.Range("C13").ClearContents()
$ws.Range("D13").ClearContents()
$ws.Range("F13").ClearContents()

$ws.Range("C14").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("F14").ClearContents()

# Match the saved selection/active-sheet state: "week 7" active, F7 selected.
$ws.Range("F7").Select()
$ws.Activate()
